$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1023.8823
$ws.Range("I43").Value = 1266.6666
$ws.Range("J43").Value = 891.4545000000001
$ws.Range("K43").Value = 1266.6666
$ws.Range("L43").Value = 891.4545000000001
$ws.Range("M43").Value = -1197.6666
$ws.Range("N43").Value = -1029.4545

$ws.Range("H82").Value = 934.5714
$ws.Range("I82").Value = 934.5714
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 2803.7142
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -2397.7142
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 934.5714
$ws.Range("I85").Value = 934.5714
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 2803.7142
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -1399.7142
$ws.Range("N85").ClearContents()

$ws.Range("H88").Value = 14001.167
$ws.Range("J88").Value = 16001
$ws.Range("L88").Value = 16001
$ws.Range("N88").Value = -16813

$ws.Range("H91").Value = 14001.167
$ws.Range("J91").Value = 16001
$ws.Range("L91").Value = 16001
$ws.Range("N91").Value = -18809

$ws.Range("H125").Value = 1087.3334
$ws.Range("J125").Value = 1087.3334
$ws.Range("L125").Value = 9786.000599999999
$ws.Range("N125").Value = -14706.0006

$ws.Range("H129").Value = 1098.5
$ws.Range("J129").Value = 2000
$ws.Range("L129").Value = 6000
$ws.Range("N129").Value = -16000

$ws.Range("H132").Value = 4229.7617
$ws.Range("I132").Value = 1855.2433
$ws.Range("J132").Value = 21801.2
$ws.Range("K132").Value = 5565.7299
$ws.Range("L132").Value = 65403.60000000001
$ws.Range("M132").Value = -3035.7299
$ws.Range("N132").Value = -70463.60000000001

$ws.Range("H138").Value = 2249679.2
$ws.Range("I138").Value = 1127.3549
$ws.Range("J138").Value = 3451491.5
$ws.Range("K138").Value = 3382.0647
$ws.Range("L138").Value = 10354474.5
$ws.Range("M138").Value = 1757.9353
$ws.Range("N138").Value = -10364754.5

$ws.Range("H141").Value = 1367.8485
$ws.Range("I141").Value = 773.16
$ws.Range("J141").Value = 3226.25
$ws.Range("K141").Value = 2319.48
$ws.Range("L141").Value = 9678.75
$ws.Range("M141").Value = 2860.52
$ws.Range("N141").Value = -20038.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1736.3226
$ws.Range("I45").Value = 1459.6364
$ws.Range("J45").Value = 2412.6667
$ws.Range("K45").Value = 1459.6364
$ws.Range("L45").Value = 2412.6667
$ws.Range("M45").Value = -1082.6364
$ws.Range("N45").Value = -3166.6667

$ws.Range("H61").Value = 1310
$ws.Range("I61").Value = 1176.4445
$ws.Range("J61").Value = 1739.2858
$ws.Range("K61").Value = 1176.4445
$ws.Range("L61").Value = 1739.2858
$ws.Range("M61").Value = -964.4445000000001
$ws.Range("N61").Value = -2163.2858

$ws.Range("H74").Value = 20990.137
$ws.Range("I74").Value = 28062.19
$ws.Range("J74").Value = 2299.7144
$ws.Range("K74").Value = 28062.19
$ws.Range("L74").Value = 2299.7144
$ws.Range("M74").Value = -27188.19
$ws.Range("N74").Value = -4047.7144

$ws.Range("H77").Value = 20990.137
$ws.Range("I77").Value = 28062.19
$ws.Range("J77").Value = 2299.7144
$ws.Range("K77").Value = 140310.95
$ws.Range("L77").Value = 11498.572
$ws.Range("M77").Value = -135942.95
$ws.Range("N77").Value = -20234.572

$ws.Range("H102").Value = 1246.25
$ws.Range("I102").Value = 1246.6666
$ws.Range("J102").Value = 1245
$ws.Range("K102").Value = 1246.6666
$ws.Range("L102").Value = 1245
$ws.Range("M102").Value = 375.3334
$ws.Range("N102").Value = -4489

$ws.Range("H123").Value = 20833.334
$ws.Range("J123").Value = 20833.334
$ws.Range("L123").Value = 20833.334
$ws.Range("N123").Value = -30633.334

$ws.Range("H132").Value = 1088.591
$ws.Range("I132").Value = 1088.591
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3265.773
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -735.7729999999997
$ws.Range("N132").ClearContents()

$ws.Range("H136").Value = 1310
$ws.Range("I136").Value = 1176.4445
$ws.Range("J136").Value = 1739.2858
$ws.Range("K136").Value = 3529.3335
$ws.Range("L136").Value = 5217.857400000001
$ws.Range("M136").Value = -979.3335000000002
$ws.Range("N136").Value = -10317.8574

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 609510.8
$ws.Range("I134").Value = 1055575.5
$ws.Range("J134").Value = 4137.1787
$ws.Range("K134").Value = 3166726.5
$ws.Range("L134").Value = 12411.5361
$ws.Range("M134").Value = -3164191.5
$ws.Range("N134").Value = -17481.5361

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1954.2623
$ws.Range("J31").Value = 2429.4595
$ws.Range("L31").Value = 2429.4595
$ws.Range("N31").Value = -3019.4595

$ws.Range("H34").Value = 1954.2623
$ws.Range("J34").Value = 2429.4595
$ws.Range("L34").Value = 2429.4595
$ws.Range("N34").Value = -2833.4595

$ws.Range("H58").Value = 5340.2
$ws.Range("I58").Value = 6583.737
$ws.Range("J58").Value = 1402.3334
$ws.Range("K58").Value = 6583.737
$ws.Range("L58").Value = 1402.3334
$ws.Range("M58").Value = -6380.737
$ws.Range("N58").Value = -1808.3334

$ws.Range("H132").Value = 608566.2
$ws.Range("I132").Value = 1309.7755
$ws.Range("J132").Value = 3088196.5
$ws.Range("K132").Value = 3929.3265
$ws.Range("L132").Value = 9264589.5
$ws.Range("M132").Value = -1399.3265
$ws.Range("N132").Value = -9269649.5

$ws.Range("H134").Value = 1438.7241
$ws.Range("I134").Value = 1382.0303
$ws.Range("K134").Value = 4146.090899999999
$ws.Range("M134").Value = -1611.090899999999

$ws.Range("H136").Value = 5340.2
$ws.Range("I136").Value = 6583.737
$ws.Range("J136").Value = 1402.3334
$ws.Range("K136").Value = 19751.211
$ws.Range("L136").Value = 4207.0002
$ws.Range("M136").Value = -17201.211
$ws.Range("N136").Value = -9307.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2591.5557
$ws.Range("I70").Value = 1332
$ws.Range("J70").Value = 7000
$ws.Range("K70").Value = 3996
$ws.Range("L70").Value = 21000
$ws.Range("M70").Value = -3681
$ws.Range("N70").Value = -21630

$ws.Range("H73").Value = 2591.5557
$ws.Range("I73").Value = 1332
$ws.Range("J73").Value = 7000
$ws.Range("K73").Value = 3996
$ws.Range("L73").Value = 21000
$ws.Range("M73").Value = -2904
$ws.Range("N73").Value = -23184

$ws.Range("H76").Value = 3000
$ws.Range("I76").Value = 3000
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 9000
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -8617
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 3000
$ws.Range("I79").Value = 3000
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 9000
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -7674
$ws.Range("N79").ClearContents()

$ws.Range("H82").Value = 3767.3333
$ws.Range("I82").Value = 1013
$ws.Range("J82").Value = 5144.5
$ws.Range("K82").Value = 3039
$ws.Range("L82").Value = 15433.5
$ws.Range("M82").Value = -2633
$ws.Range("N82").Value = -16245.5

$ws.Range("H85").Value = 3767.3333
$ws.Range("I85").Value = 1013
$ws.Range("J85").Value = 5144.5
$ws.Range("K85").Value = 3039
$ws.Range("L85").Value = 15433.5
$ws.Range("M85").Value = -1635
$ws.Range("N85").Value = -18241.5

$ws.Range("H88").Value = 5950
$ws.Range("J88").Value = 5950
$ws.Range("L88").Value = 17850
$ws.Range("N88").Value = -18706

$ws.Range("H91").Value = 5950
$ws.Range("J91").Value = 5950
$ws.Range("L91").Value = 17850
$ws.Range("N91").Value = -20814

$ws.Range("H113").Value = 2525764.5
$ws.Range("I113").Value = 6061154.5
$ws.Range("J113").Value = 485.85715
$ws.Range("K113").Value = 18183463.5
$ws.Range("L113").Value = 1457.57145
$ws.Range("M113").Value = -18181293.5
$ws.Range("N113").Value = -5797.571449999999

$ws.Range("H131").Value = 919.85
$ws.Range("I131").Value = 815
$ws.Range("J131").Value = 921.9897999999999
$ws.Range("K131").Value = 2445
$ws.Range("L131").Value = 2765.9694
$ws.Range("M131").Value = 2595
$ws.Range("N131").Value = -12845.9694

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2613.75
$ws.Range("I126").Value = 3477.75
$ws.Range("J126").Value = 2181.75
$ws.Range("K126").Value = 10433.25
$ws.Range("L126").Value = 6545.25
$ws.Range("M126").Value = -7963.25
$ws.Range("N126").Value = -11485.25

$ws.Range("H132").Value = 1641799.4
$ws.Range("I132").Value = 2155.805
$ws.Range("J132").Value = 5003069
$ws.Range("K132").Value = 6467.414999999999
$ws.Range("L132").Value = 15009207
$ws.Range("M132").Value = -3937.414999999999
$ws.Range("N132").Value = -15014267

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1355.431
$ws.Range("I136").Value = 997.8946999999999
$ws.Range("K136").Value = 2993.6841
$ws.Range("M136").Value = -443.6840999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1007.3333
$ws.Range("I126").Value = 946.36365
$ws.Range("J126").Value = 1175
$ws.Range("K126").Value = 2839.09095
$ws.Range("L126").Value = 3525
$ws.Range("M126").Value = -369.0909499999998
$ws.Range("N126").Value = -8465

$ws.Range("H136").Value = 2075.7656
$ws.Range("I136").Value = 1829.2391
$ws.Range("J136").Value = 2705.7778
$ws.Range("K136").Value = 5487.7173
$ws.Range("L136").Value = 8117.3334
$ws.Range("M136").Value = -2937.7173
$ws.Range("N136").Value = -13217.3334
